# Update countries & provincias Spain
# Daily COVID data refresh: "Pais" sheet numbers updated, timestamp bumped,
# and two country pairs swap rank (and therefore row position) because one
# of each pair grew past its former neighbour.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp (A1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 00:52"

# --- Estados Unidos (row 4) -------------------------------------------------
$ws.Range("B4").Value = 558523
$ws.Range("C4").Value = 25644
$ws.Range("E4").Value = 504544
$ws.Range("G4").Value = 1416
$ws.Range("H4").Value = 21993

# --- Noruega (row 31) --------------------------------------------------------
$ws.Range("B31").Value = 6525
$ws.Range("C31").Value = 116
$ws.Range("E31").Value = 6365

# --- Chequia (row 35) ---------------------------------------------------------
$ws.Range("B35").Value = 5991
$ws.Range("C35").Value = 89
$ws.Range("E35").Value = 5389

# --- Hong Kong / Bosnia y Herzegovina swap rank (rows 73 & 74) --------------
# Bosnia y Herzegovina overtakes Hong Kong (1009 vs 1005 total cases), so it
# now occupies row 73 while Hong Kong (unchanged numbers) drops to row 74.
$ws.Range("A73").Value = "Bosnia y Herzegovina"
$ws.Range("B73").Value = 1009
$ws.Range("C73").Value = 63
$ws.Range("D73").Value = 193
$ws.Range("E73").Value = 777
$ws.Range("F73").Value = 4
$ws.Range("G73").Value = 2
$ws.Range("H73").Value = 39

$ws.Range("A74").Value = "Hong Kong"
$ws.Range("B74").Value = 1005
$ws.Range("C74").Value = 4
$ws.Range("D74").Value = 360
$ws.Range("E74").Value = 641
$ws.Range("F74").Value = 13
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 4

# --- Reunion / Jordania swap rank (rows 99 & 100) ---------------------------
# Reunion overtakes Jordania (389 vs 389 total cases but higher in other
# columns) so it now occupies row 99 while Jordania (unchanged numbers)
# drops to row 100.
$ws.Range("A99").Value = "Reunion"
$ws.Range("B99").Value = 389
$ws.Range("C99").Value = 1
$ws.Range("D99").Value = 40
$ws.Range("E99").Value = 349
$ws.Range("F99").Value = 3
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0

$ws.Range("A100").Value = "Jordania"
$ws.Range("B100").Value = 389
$ws.Range("C100").Value = 8
$ws.Range("D100").Value = 201
$ws.Range("E100").Value = 181
$ws.Range("F100").Value = 5
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 7

# --- Venezuela (row 121) -----------------------------------------------------
$ws.Range("B121").Value = 181
$ws.Range("C121").Value = 6
$ws.Range("E121").Value = 79

# --- Monaco (row 134) --------------------------------------------------------
$ws.Range("B134").Value = 93
$ws.Range("C134").Value = 1
$ws.Range("D134").Value = 6
$ws.Range("F134").Value = 5

# --- Dominica (row 179) ------------------------------------------------------
$ws.Range("D179").Value = 8
$ws.Range("E179").Value = 8
